$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1416.7778
$ws.Range("I19").Value = 1350
$ws.Range("K19").Value = 1350
$ws.Range("M19").Value = -1175
$ws.Range("H33").Value = 237.58824
$ws.Range("I33").Value = 211.33333
$ws.Range("J33").Value = 434.5
$ws.Range("K33").Value = 211.33333
$ws.Range("L33").Value = 434.5
$ws.Range("M33").Value = 17.66667000000001
$ws.Range("N33").Value = -892.5
$ws.Range("H43").Value = 6766.35
$ws.Range("J43").Value = 6763.143
$ws.Range("L43").Value = 6763.143
$ws.Range("N43").Value = -6901.143
$ws.Range("H64").Value = 7583.9443
$ws.Range("I64").Value = 4342
$ws.Range("J64").Value = 12678.429
$ws.Range("K64").Value = 4342
$ws.Range("L64").Value = 12678.429
$ws.Range("M64").Value = -4094
$ws.Range("N64").Value = -13174.429
$ws.Range("H67").Value = 7583.9443
$ws.Range("I67").Value = 4342
$ws.Range("J67").Value = 12678.429
$ws.Range("K67").Value = 4342
$ws.Range("L67").Value = 12678.429
$ws.Range("M67").Value = -3484
$ws.Range("N67").Value = -14394.429
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("H132").Value = 2041.0555
$ws.Range("I132").Value = 1799.1538
$ws.Range("J132").Value = 2670
$ws.Range("K132").Value = 5397.4614
$ws.Range("L132").Value = 8010
$ws.Range("M132").Value = -2867.4614
$ws.Range("N132").Value = -13070
$ws.Range("H138").Value = 2055.9375
$ws.Range("I138").Value = 1607.3077
$ws.Range("K138").Value = 4821.9231
$ws.Range("M138").Value = 318.0769
$ws.Range("H141").Value = 721.0952
$ws.Range("I141").Value = 721.0952
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2163.2856
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3016.7144
$ws.Range("M87").ClearContents()
$ws.Range("M90").ClearContents()
$ws.Range("N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2038.1154
$ws.Range("I2").Value = 1360
$ws.Range("J2").Value = 2962.818
$ws.Range("K2").Value = 1360
$ws.Range("L2").Value = 2962.818
$ws.Range("M2").Value = -1247
$ws.Range("N2").Value = -3188.818
$ws.Range("H74").Value = 1052.4117
$ws.Range("I74").Value = 1052.4117
$ws.Range("K74").Value = 1052.4117
$ws.Range("M74").Value = -178.4117000000001
$ws.Range("H77").Value = 1052.4117
$ws.Range("I77").Value = 1052.4117
$ws.Range("K77").Value = 5262.058500000001
$ws.Range("M77").Value = -894.058500000001
$ws.Range("H116").Value = 2038.1154
$ws.Range("I116").Value = 1360
$ws.Range("J116").Value = 2962.818
$ws.Range("K116").Value = 1360
$ws.Range("L116").Value = 2962.818
$ws.Range("M116").Value = 934
$ws.Range("N116").Value = -7550.818
$ws.Range("H132").Value = 1715
$ws.Range("I132").Value = 1019.75
$ws.Range("K132").Value = 3059.25
$ws.Range("M132").Value = -529.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2038.1154
$ws.Range("I3").Value = 1360
$ws.Range("J3").Value = 2962.818
$ws.Range("K3").Value = 1360
$ws.Range("L3").Value = 2962.818
$ws.Range("M3").Value = -1246
$ws.Range("N3").Value = -3190.818
$ws.Range("H74").Value = 52191.75
$ws.Range("J74").Value = 52191.75
$ws.Range("L74").Value = 52191.75
$ws.Range("N74").Value = -54063.75
$ws.Range("H77").Value = 52191.75
$ws.Range("J77").Value = 52191.75
$ws.Range("L77").Value = 156575.25
$ws.Range("N77").Value = -165935.25
$ws.Range("H105").Value = 3151.4443
$ws.Range("I105").Value = 2826.6875
$ws.Range("J105").Value = 5749.5
$ws.Range("K105").Value = 2826.6875
$ws.Range("L105").Value = 5749.5
$ws.Range("M105").Value = -1079.6875
$ws.Range("N105").Value = -9243.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1539.0667
$ws.Range("I31").Value = 1506.6154
$ws.Range("J31").Value = 1750
$ws.Range("K31").Value = 1506.6154
$ws.Range("L31").Value = 1750
$ws.Range("M31").Value = -1211.6154
$ws.Range("N31").Value = -2340
$ws.Range("H34").Value = 1539.0667
$ws.Range("I34").Value = 1506.6154
$ws.Range("J34").Value = 1750
$ws.Range("K34").Value = 1506.6154
$ws.Range("L34").Value = 1750
$ws.Range("M34").Value = -1304.6154
$ws.Range("N34").Value = -2154
$ws.Range("H58").Value = 1063.6666
$ws.Range("I58").Value = 726.45
$ws.Range("K58").Value = 726.45
$ws.Range("M58").Value = -523.45
$ws.Range("H60").Value = 19666.666
$ws.Range("I60").Value = 19666.666
$ws.Range("K60").Value = 19666.666
$ws.Range("M60").Value = -19155.666
$ws.Range("H99").Value = 7596.8823
$ws.Range("I99").Value = 7476.467
$ws.Range("K99").Value = 7476.467
$ws.Range("M99").Value = -5978.467
$ws.Range("H126").Value = 7596.8823
$ws.Range("I126").Value = 7476.467
$ws.Range("K126").Value = 22429.401
$ws.Range("M126").Value = -19959.401
$ws.Range("H136").Value = 1063.6666
$ws.Range("I136").Value = 726.45
$ws.Range("K136").Value = 2179.35
$ws.Range("M136").Value = 370.6499999999996

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 8098.4
$ws.Range("J80").Value = 9499.5
$ws.Range("L80").Value = 28498.5
$ws.Range("N80").Value = -30370.5
$ws.Range("H83").Value = 8098.4
$ws.Range("J83").Value = 9499.5
$ws.Range("L83").Value = 85495.5
$ws.Range("N83").Value = -94855.5
$ws.Range("H107").Value = 1322.3334
$ws.Range("J107").Value = 1322.3334
$ws.Range("L107").Value = 3967.0002
$ws.Range("N107").Value = -7807.0002
$ws.Range("H113").Value = 846
$ws.Range("I113").Value = 774
$ws.Range("K113").Value = 2322
$ws.Range("M113").Value = -152
$ws.Range("H131").Value = 924.4
$ws.Range("I131").Value = 826
$ws.Range("K131").Value = 2478
$ws.Range("M131").Value = 2562

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5138.8
$ws.Range("I40").Value = 5138.8
$ws.Range("K40").Value = 5138.8
$ws.Range("M40").Value = -5002.8
$ws.Range("H46").Value = 4793.1035
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 4400
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 4400
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -4776

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2499.8333
$ws.Range("I62").Value = 2399.8
$ws.Range("K62").Value = 2399.8
$ws.Range("M62").Value = -1775.8
$ws.Range("H65").Value = 2499.8333
$ws.Range("I65").Value = 2399.8
$ws.Range("K65").Value = 11999
$ws.Range("M65").Value = -8879
$ws.Range("H132").Value = 1064.2439
$ws.Range("I132").Value = 1099.8379
$ws.Range("J132").Value = 735
$ws.Range("K132").Value = 3299.5137
$ws.Range("L132").Value = 2205
$ws.Range("M132").Value = -769.5137
$ws.Range("N132").Value = -7265
